$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B. This shifts the original Execute column (B) to
# C, and the original InvocationCount column (C, whose data cells carried the
# quotePrefix style) to D -- preserving that style assignment on the shifted
# cells exactly like Excel does for a real column insert.
$ws.Columns("B").Insert()

# --- Header row ---
$ws.Range("B1").Value = "TestDescription"

# --- New TestDescription column (added to shared strings before "login") ---
$ws.Range("B2").Value = "This is a test to validate login in to OrangeHrm"
$ws.Range("B3").Value = "This is a copy "

# --- Data rows: TestName becomes the shared "login" value ---
$ws.Range("A2").Value = "login"
$ws.Range("A3").Value = "login"

# --- InvocationCount (shifted to D) becomes numeric 2 / 1, keeping the
# original quotePrefix-style formatting that lived on those cells. ---
$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("D2").Value = 2
$ws.Range("F2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("F2").Clear()

$ws.Range("D3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("D3").Value = 1
$ws.Range("F3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("F3").Clear()

$ws.Range("D2:D3").NumberFormat = "@"

# --- Column B width to match column A's custom width ---
$ws.Columns("B").ColumnWidth = 19.6640625

# --- Selection / view state ---
$ws.Range("F12").Select() | Out-Null

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1
